$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36: new log entry (2022-07-01, 21:00 - 22:45) ---
$ws.Range("A35:G35").Copy()
$ws.Range("A36:G36").PasteSpecial(-4122)
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = 44743
$ws.Range("C36").Value = 0.875
$ws.Range("D36").Value = 0.94791666666666663
$ws.Range("E36").Formula = "=D36-C36"
$ws.Range("F36").Value = "Code"
$ws.Range("G36").Value = "Presentation content for Intro, FCN, Unet and PSPNet completed"

# --- Row 37: new log entry (2022-07-02, 08:00 - 09:00) ---
$ws.Range("A35:G35").Copy()
$ws.Range("A37:G37").PasteSpecial(-4122)
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = 44744
$ws.Range("C37").Value = 0.33333333333333331
$ws.Range("D37").Value = 0.375
$ws.Range("E37").Formula = "=D37-C37"
$ws.Range("F37").Value = "Code"
$ws.Range("G37").Value = "1. Test video starter nb for costa rica and paris challenge`n2. paris_challenge.mov = (2160, 3840, 3) at 60fps`n3. costa rica avi (1080, 1900, 3) at 30fps`n4. need to resize images to (640, 360) for test video scenes matching"
$ws.Rows.Item(37).RowHeight = 75

$excel.CutCopyMode = 0

# --- Total row (row 38) formula auto-recalculates since it sums E2:E37 ---

# --- Update selection to match the last-edited cell ---
$ws.Range("D38").Select()
